# Automatic update of files.
#
# Upstream re-synced this export and the two "Knärot" (Goodyera repens)
# observation rows (4 and 5) ended up re-ordered: the fields that are
# specific to each individual observation (Id, Antal, the Ost/Nord
# coordinate pair and the Starttid/Sluttid timestamps) got exchanged
# between row 4 and row 5, while every shared attribute (species info,
# locality, dates, observer, ...) is identical between the two rows and
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric columns: Id (A), Ost (Q), Nord (R) -----------------------
$numericCols = "A", "Q", "R"
foreach ($col in $numericCols) {
    $cell4 = $ws.Range("$col`4")
    $cell5 = $ws.Range("$col`5")

    $val4 = $cell4.Value2
    $val5 = $cell5.Value2

    $cell4.Value = $val5
    $cell5.Value = $val4
}

# --- Starttid / Sluttid (Z, AB): plain text "hh:mm" values ------------
$timeCols = "Z", "AB"
foreach ($col in $timeCols) {
    $cell4 = $ws.Range("$col`4")
    $cell5 = $ws.Range("$col`5")

    $val4 = [string]$cell4.Value2
    $val5 = [string]$cell5.Value2

    $cell4.Value = $val5
    $cell5.Value = $val4
}

# --- Antal (I): text cell, value "35" on row 5 / empty on row 4 -------
# Read as text explicitly (Value2 would coerce "35" to a Double) and
# write back with a leading apostrophe so the numeric-looking string is
# stored as text again (matches the original inlineStr "35"), instead of
# being reinterpreted as a number.
$i4 = $ws.Range("I4")
$i5 = $ws.Range("I5")

$i4Text = $i4.Text
$i5Text = $i5.Text

$i4.Value = "'" + $i5Text
$i5.Value = "'" + $i4Text
